$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay as literal text
# (Excel would otherwise coerce them to Double and lose the exact display,
# e.g. "1.0000" -> 1 or "6.440" -> 6.44). Force Text format on those first.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D11", "D14", "D15", "D16", "D18", "D20", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D36", "D37", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values exactly as scraped
$ws.Range('D2').Value = '29.387.26'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.847.67'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '240.17'
$ws.Range('D6').Value = '0.6303'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.07596'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '0.2931'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('D11').Value = '0.07743'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '1.834.85'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Value = '0.00001079'
$ws.Range('E14').Value = '  +7.50%  '
$ws.Range('D15').Value = '0.6786'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '83.64'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').Value = '2.091.72'
$ws.Range('E17').Value = '  -7.62%  '
$ws.Range('D18').Value = '6.151'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').Value = '29.407.08'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '229.09'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D23').Value = '7.427'
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '157.11'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '0.1388'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('D27').Value = '8.383'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '17.62'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').Value = '1.310'
$ws.Range('E29').Value = '  +4.71%  '
$ws.Range('D30').Value = '1.463'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').Value = '4.106'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').Value = '4.033'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').Value = '1.848'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D36').Value = '0.7098'
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').Value = '2.582'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '1.238.20'
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('D39').Value = '2.772'
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('D41').Value = '6.449'
$ws.Range('E41').Value = '  +3.82%  '
$ws.Range('D42').Value = '0.9085'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').Value = '1.0000'
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = '101.54'
$ws.Range('D45').Value = '66.13'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').Value = '0.00000000123'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D47').Value = '7.190'
$ws.Range('E47').Value = '  +1.79%  '
$ws.Range('D48').Value = '0.4013'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.680'
$ws.Range('E49').Value = '  -1.65%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.961'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('E51').Value = '  -0.78%  '
